# Update automàtic: dades i banners [2026-02-20 09:45]
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dades_Període")

# --- Row 2: period/time & timestamp text cells ---
$ws.Range("E2").Value  = "09:00 - 09:30"
$ws.Range("H2").Value  = "2026-02-20 09:45:28"
$ws.Range("Q2").Value  = "09:00 - 09:30"
$ws.Range("X2").Value  = "09:00 - 09:30"
$ws.Range("AI2").Value = "09:00 - 09:30"

# --- Row 2: numeric-looking readings -> keep them as text values      ---
# (leading apostrophe forces Excel to store them as text, same as the  ---
#  source workbook, instead of auto-converting to a number)            ---
$ws.Range("M2").Value  = "'131"
$ws.Range("N2").Value  = "'58"
$ws.Range("O2").Value  = "'1024.1"
$ws.Range("R2").Value  = "'345"
$ws.Range("S2").Value  = "'10.9"
$ws.Range("T2").Value  = "'10.5"
$ws.Range("U2").Value  = "'11.6"
$ws.Range("V2").Value  = "'1.8"
$ws.Range("W2").Value  = "'4.7"
$ws.Range("Y2").Value  = "'10.9"
$ws.Range("Z2").Value  = "'11.6"
$ws.Range("AA2").Value = "'10.5"
$ws.Range("AB2").Value = "'58"
$ws.Range("AD2").Value = "'1.8"
$ws.Range("AE2").Value = "'131"
$ws.Range("AF2").Value = "'4.7"
$ws.Range("AG2").Value = "'1024.1"
$ws.Range("AH2").Value = "'345"
$ws.Range("AJ2").Value = "'10.9"
$ws.Range("AK2").Value = "'11.6"
$ws.Range("AL2").Value = "'10.5"
$ws.Range("AM2").Value = "'58"
$ws.Range("AO2").Value = "'1.8"
$ws.Range("AP2").Value = "'131"
$ws.Range("AQ2").Value = "'4.7"
$ws.Range("AR2").Value = "'1024.1"
$ws.Range("AS2").Value = "'345"

# --- Rows 3-6: DATA_EXTRACCIO timestamp refresh ---
$ws.Range("H3").Value = "2026-02-20 09:45:30"
$ws.Range("H4").Value = "2026-02-20 09:45:30"
$ws.Range("H5").Value = "2026-02-20 09:45:30"
$ws.Range("H6").Value = "2026-02-20 09:45:30"
